$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B/C/D/E values per row (G = B+C+D+E is recomputed automatically by Excel
# if there were a formula; here G is a literal value too, so we set it explicitly).
$data = @{
    2  = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    3  = @(3.182878228561681, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    4  = @(3.182878228561681, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    5  = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    6  = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    7  = @(1.505614041169197, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    8  = @(0.7287194209349384, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 2.27892381503245)
    9  = @(3.182878228561681, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    10 = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    11 = @(0.06328177979961902, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.613486173897131)
    12 = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    13 = @(0.7287194209349384, 0.004309184025731883, 3.082599426703578, 0.4998867070740569, 4.315514738738305)
    14 = @(3.182878228561681, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    15 = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
